$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Price" column (D) stores numbers as plain text using "." as both
# a thousands separator and (for sub-$1 coins) a decimal point, e.g.
# "27.894.69" or "1.001". Briefly switch affected cells to text format
# before writing so Excel does not silently reinterpret them as numeric
# values and mangle them (e.g. "1.001" -> 1.0009999999999999), then
# restore "General" so the cell formatting matches the original file.
function Set-TextValue($addr, $value) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.NumberFormat = "General"
}

# Updated coin prices / 1h volume changes from the latest GitHub Actions run.
Set-TextValue "D2" '27.894.69'
$ws.Range("E2").Value = '  +1.31%  '
Set-TextValue "D3" '1.753.11'
$ws.Range("E3").Value = '  -1.16%  '
$ws.Range("E4").Value = '  -0.70%  '
Set-TextValue "D5" '335.92'
$ws.Range("E5").Value = '  -0.45%  '
Set-TextValue "D6" '1.001'
$ws.Range("E6").Value = '  -0.62%  '
Set-TextValue "D7" '0.3821'
$ws.Range("E7").Value = '  -1.33%  '
Set-TextValue "D8" '0.3391'
$ws.Range("E8").Value = '  -1.06%  '
Set-TextValue "D9" '44.43'
$ws.Range("E9").Value = '  -5.80%  '
Set-TextValue "D10" '1.110'
$ws.Range("E10").Value = '  -3.73%  '
Set-TextValue "D11" '0.07216'
$ws.Range("E11").Value = '  -3.66%  '
$ws.Range("E12").Value = '  -0.51%  '
Set-TextValue "D13" '22.45'
$ws.Range("E13").Value = '  -0.19%  '
Set-TextValue "D14" '6.149'
$ws.Range("E14").Value = '  -4.16%  '
Set-TextValue "D15" '7.132'
$ws.Range("E15").Value = '  -0.20%  '
Set-TextValue "D16" '1.757.57'
$ws.Range("E16").Value = '  -1.26%  '
Set-TextValue "D17" '0.00001056'
$ws.Range("E17").Value = '  -2.23%  '
Set-TextValue "D18" '0.06608'
$ws.Range("E18").Value = '  -1.64%  '
Set-TextValue "D19" '79.04'
$ws.Range("E19").Value = '  -5.14%  '
$ws.Range("E20").Value = '  -0.44%  '
$ws.Range("E21").Value = '  -5.23%  '
Set-TextValue "D22" '6.220'
$ws.Range("E22").Value = '  -4.12%  '
Set-TextValue "D23" '27.926.88'
$ws.Range("E23").Value = '  +1.28%  '
Set-TextValue "D24" '11.60'
$ws.Range("E24").Value = '  -5.22%  '
Set-TextValue "D25" '2.382'
$ws.Range("E25").Value = '  -0.16%  '
Set-TextValue "D26" '152.66'
$ws.Range("E26").Value = '  -1.16%  '
Set-TextValue "D27" '19.85'
$ws.Range("E27").Value = '  -5.17%  '
Set-TextValue "D28" '2.315'
$ws.Range("E28").Value = '  -6.80%  '
Set-TextValue "D29" '1.958.66'
$ws.Range("E29").Value = '  -1.12%  '
Set-TextValue "D30" '1.272'
$ws.Range("E30").Value = '  -11.49%  '
Set-TextValue "D31" '131.78'
$ws.Range("E31").Value = '  -3.19%  '
Set-TextValue "D32" '4.017'
$ws.Range("E32").Value = '  +0.95%  '
Set-TextValue "D33" '5.811'
$ws.Range("E33").Value = '  -6.43%  '
Set-TextValue "D34" '0.08805'
$ws.Range("E34").Value = '  -1.49%  '
Set-TextValue "D35" '12.19'
$ws.Range("E35").Value = '  -5.60%  '
Set-TextValue "D36" '0.6591'
$ws.Range("E36").Value = '  -3.52%  '
Set-TextValue "D37" '0.06171'
$ws.Range("E37").Value = '  -3.89%  '
$ws.Range("B38").Value = 'VeChain'
$ws.Range("C38").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
Set-TextValue "D38" '0.02281'
$ws.Range("E38").Value = '  -7.05%  '
$ws.Range("B39").Value = 'InternetComputer(DFINITY)'
$ws.Range("C39").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
Set-TextValue "D39" '5.119'
$ws.Range("E39").Value = '  -5.83%  '
Set-TextValue "D40" '1.514'
$ws.Range("E40").Value = '  -1.18%  '
Set-TextValue "D41" '0.2105'
$ws.Range("E41").Value = '  -4.90%  '
Set-TextValue "D42" '1.208'
$ws.Range("E42").Value = '  -3.48%  '
Set-TextValue "D43" '7.970'
$ws.Range("E43").Value = '  -5.67%  '
$ws.Range("E44").Value = '  -0.45%  '
Set-TextValue "D45" '13.75'
$ws.Range("E45").Value = '  -4.65%  '
Set-TextValue "D46" '3.824'
$ws.Range("E46").Value = '  -0.93%  '
Set-TextValue "D47" '0.6047'
$ws.Range("E47").Value = '  -4.66%  '
Set-TextValue "D48" '125.93'
$ws.Range("E48").Value = '  -4.97%  '
Set-TextValue "D49" '2.000'
$ws.Range("E49").Value = '  -5.75%  '
Set-TextValue "D50" '1.123'
$ws.Range("E50").Value = '  +5.49%  '
Set-TextValue "D51" '1.166'
$ws.Range("E51").Value = '  +0.46%  '
